$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column D ("Tipo"), shifting Tipo to E
$ws.Columns.Item(4).Insert()

# Header for the newly inserted column D ("MAE"), matching the other header cells' formatting
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").Borders.Weight = 2

# Updated metric values for row 2 (AVG_TIME_TO_NEU)
$ws.Range("B2").Value = 0.1058112132799173
$ws.Range("C2").Value = 0.9980670797160867
$ws.Range("D2").Value = 0.2321132172107814

# Updated metric values for row 3 (AVG_TIME_TO_LEA)
$ws.Range("B3").Value = 0.2229918473948952
$ws.Range("C3").Value = 0.9835245499612149
$ws.Range("D3").Value = 0.3651187035608621
